# Apply updates to the distribution table on Sheet1 ("Final version of presentation")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("U2").Value = 0.23
$ws.Range("V2").Value = 0.35
$ws.Range("W2").Value = 0.15
$ws.Range("X2").Value = 0.12
$ws.Range("Y2").Value = 0.05
$ws.Range("Z2").Value = 0.03
$ws.Range("AA2").Value = 0.02
$ws.Range("AB2").Value = 0.01
$ws.Range("AC2").Value = 0.01
$ws.Range("AD2").Value = 0.01

# Row 3
$ws.Range("V3").Value = 0.35
$ws.Range("W3").Value = 0.15
$ws.Range("X3").Value = 0.12
$ws.Range("Y3").Value = 0.05
$ws.Range("Z3").Value = 0.03
$ws.Range("AA3").Value = 0.02
$ws.Range("AB3").Value = 0.01
$ws.Range("AC3").Value = 0.01
$ws.Range("AD3").Value = 0.01

# Row 4
$ws.Range("V4").Value = 0.06

# Row 5
$ws.Range("U5").Value = 0.23
$ws.Range("V5").Value = 0.35
$ws.Range("W5").Value = 0.15
$ws.Range("X5").Value = 0.12
$ws.Range("Y5").Value = 0.05
$ws.Range("Z5").Value = 0.03
$ws.Range("AA5").Value = 0.02
$ws.Range("AB5").Value = 0.01
$ws.Range("AC5").Value = 0.01
